$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "2.5% (1035)",
    "17% (110)",
    "39% (1091)",
    "85% (943)",
    "39% (800)",
    "38% (1090)",
    "36% (1028)",
    "62% (144)",
    "55% (56)",
    "58% (1089)",
    "15% (537)",
    "16% (1085)",
    "0.19% (1030)",
    "9.7% (1091)",
    "3.6% (1091)",
    "0.87% (923)",
    "0.34% (290)",
    "2.2% (458)",
    "6.5% (1078)",
    "0% (1088)",
    "2.5% (1091)",
    "0% (1091)",
    "0.7% (718)",
    "0.092% (1091)",
    "0% (390)",
    "0.29% (345)",
    "0% (1090)"
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

# Row 29 (F29) shares the same original value as row 23 ("0 (1091)") and
# must be updated to match the corresponding new value ("0% (1091)").
$ws.Cells.Item(29, 6).Value = "0% (1091)"
